$d = $word.ActiveDocument

# 1. Merge title run + trailing-space run into a single run with trailing space.
$d.Content.Find.Execute(
    "Использование средств автоматической генерации документации и рефакторинга программного кода ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Использование средств автоматической генерации документации и рефакторинга программного кода ",
    2) | Out-Null

# 2. "Вариант 3" paragraph -> remove text/runs, drop centering, add en-US lang on mark.
$rng = $d.Content
$rng.Find.Execute("Вариант 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p = $rng.Paragraphs(1)
$pr = $p.Range
$pr.MoveEnd(1, -1) | Out-Null
$pr.Delete()
$p.Alignment = 0
$p.Range.LanguageID = "en-US"

# 3. The two following empty paragraphs gain an en-US language mark too.
$p2 = $p.Next()
$p2.Range.LanguageID = "en-US"
$p3 = $p2.Next()
$p3.Range.LanguageID = "en-US"
